# Hook up the spiral/mill cutting-speed lookup table on the "Cutting Speed"
# sheet: replace the old, sparsely-filled B3:E14 block with a fully
# populated B3:E22 table of (Cutting Meter, Mill Diameter, Number of teeth,
# Feed pr Tooth) rows used by the spiral calculator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old, partially-filled rows (3-14) so no stray cells remain.
$ws.Range("B3:E14").ClearContents()

# Each row: Cutting Meter, Mill Diameter, Number of teeth, Feed pr Tooth
$data = @(
  @(210, 80,  6, "0,12"),
  @(180, 80,  6, "0,12"),
  @(300, 80,  6, "0,12"),
  @(500, 80,  6, "0,12"),
  @(500, 40,  4, "0,1"),
  @(500, 20,  4, "0,1"),
  @(40,  20,  4, "0,1"),
  @(60,  20,  4, "0,1"),
  @(80,  20,  4, "0,1"),
  @(100, 20,  4, "0,1"),
  @(120, 20,  4, "0,1"),
  @(140, 20,  4, "0,1"),
  @(160, 20,  4, "0,1"),
  @(180, 20,  4, "0,1"),
  @(200, 20,  4, "0,1"),
  @(200, 100, 4, "0,1"),
  @(200, 10,  4, "0,1"),
  @(10,  10,  4, "0,1"),
  @(20,  10,  4, "0,1"),
  @(70,  10,  4, "0,1")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value2 = $row[0]
    $ws.Cells.Item($r, 3).Value2 = $row[1]
    $ws.Cells.Item($r, 4).Value2 = $row[2]
    $ws.Cells.Item($r, 5).Value2 = $row[3]
    $r = $r + 1
}
